$d = $word.ActiveDocument

# --- Table row height: 11655 twips (582.75pt) -> 11026 twips (551.3pt) ---
$table = $d.Tables.Item(1)
$row = $table.Rows.Item(1)
$row.Height = 551.3

# --- Style updates: Calibri -> Aptos, various sizes, uiPriority ---

# CLIN1HEADING (paragraph style)
$s = $d.Styles("CLIN1HEADING")
$s.Font.Name = "Aptos"
$s.Font.Size = 15

# CLIN2SUBHEADINGS (paragraph style)
$s = $d.Styles("CLIN2SUBHEADINGS")
$s.Font.Name = "Aptos"
$s.Font.Size = 10

# CLIN1HEADINGChar (character style)
$s = $d.Styles("CLIN1HEADINGChar")
$s.Font.Name = "Aptos"
$s.Font.Size = 15

# CLIN3BULLETPOINTS (paragraph style)
$s = $d.Styles("CLIN3BULLETPOINTS")
$s.Font.Name = "Aptos"
$s.Font.Size = 8

# CLIN2SUBHEADINGSChar (character style)
$s = $d.Styles("CLIN2SUBHEADINGSChar")
$s.Font.Name = "Aptos"
$s.Font.Size = 10

# CLIN4 (paragraph style)
$s = $d.Styles("CLIN4")
$s.Font.Name = "Aptos"
$s.Font.Size = 5.5

# CLIN3BULLETPOINTSChar (character style)
$s = $d.Styles("CLIN3BULLETPOINTSChar")
$s.Font.Name = "Aptos"
$s.Font.Size = 8

# CLIN4Char (character style)
$s = $d.Styles("CLIN4Char")
$s.Priority = 1
$s.Font.Name = "Aptos"
$s.Font.Size = 5.5
$s.Font.SizeBi = 6
